$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cells per the scraped symbol-list refresh.
# Numeric-looking Price/Volume columns (D/E) are forced to Text format
# before assignment so values are stored verbatim (e.g. "256.87", "-0.48%")
# instead of being auto-converted into numbers/percentages by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '256.87'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.00%'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.48%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.662'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-1.29%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05894'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.607'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.76%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8501'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-2.27%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9217'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-2.97%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1376'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-1.80%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.04289'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '9.10%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07011'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-2.07%'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-4.26%'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-1.49%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001529'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-1.14%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.006026'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.44%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.468'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.44%'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.165'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.12%'
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.194'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.14%'
$ws.Range("B19").Value = 'One'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.01030'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1,602.21%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3054'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-2.58%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.20%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.911'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '2.84%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04256'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.25%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001218'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.18%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004296'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-4.48%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.03%'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-21.33%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03782'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-1.26%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006230'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '2.16%'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-0.06%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002200'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-2.25%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01413'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '33.51%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005372'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-2.29%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.03%'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '10,454.21%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002100'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.03%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002000'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.03%'
